$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'288.72"
$ws.Range("D3").Value = "'40.39"
$ws.Range("E3").Value = "'-2.34%"
$ws.Range("D4").Value = "'5.037"
$ws.Range("E4").Value = "'-3.99%"
$ws.Range("D5").Value = "'0.07302"
$ws.Range("E5").Value = "'-5.70%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.278"
$ws.Range("E6").Value = "'-1.40%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.525"
$ws.Range("E7").Value = "'-9.92%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9186"
$ws.Range("E8").Value = "'-2.69%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1188"
$ws.Range("E9").Value = "'-4.56%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1722"
$ws.Range("E10").Value = "'-7.29%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08614"
$ws.Range("E11").Value = "'-6.73%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04165"
$ws.Range("E12").Value = "'-3.67%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1055"
$ws.Range("E13").Value = "'0.30%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001264"
$ws.Range("E14").Value = "'-1.38%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005787"
$ws.Range("E15").Value = "'-3.93%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.397"
$ws.Range("E16").Value = "'1.61%"
$ws.Range("D18").Value = "'0.3286"
$ws.Range("E18").Value = "'-1.41%"
$ws.Range("E19").Value = "'0.26%"
$ws.Range("D21").Value = "'0.2881"
$ws.Range("E21").Value = "'1.99%"
$ws.Range("E22").Value = "'-4.62%"
$ws.Range("E23").Value = "'-0.14%"
$ws.Range("D24").Value = "'0.003821"
$ws.Range("E24").Value = "'-7.31%"
$ws.Range("D25").Value = "'0.0001279"
$ws.Range("E25").Value = "'0.64%"
$ws.Range("D38").Value = "'0.02312"
$ws.Range("E38").Value = "'-9.40%"
$ws.Range("D39").Value = "'0.04974"
$ws.Range("E39").Value = "'-7.08%"
$ws.Range("D40").Value = "'0.006413"
$ws.Range("E40").Value = "'221.96%"
$ws.Range("D41").Value = "'0.007673"
$ws.Range("E41").Value = "'-1.24%"
$ws.Range("D42").Value = "'0.1274"
$ws.Range("E42").Value = "'-3.35%"
$ws.Range("D43").Value = "'0.007348"
$ws.Range("E43").Value = "'0.16%"
$ws.Range("D44").Value = "'0.007046"
$ws.Range("E44").Value = "'-14.57%"
$ws.Range("D45").Value = "'0.3123"
$ws.Range("E45").Value = "'-1.71%"
$ws.Range("D46").Value = "'0.00006431"
$ws.Range("E46").Value = "'-4.44%"
$ws.Range("E47").Value = "'-0.14%"
$ws.Range("D48").Value = "'0.2506"
$ws.Range("E48").Value = "'24.13%"
$ws.Range("E49").Value = "'-0.06%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.14%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.14%"
